$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("empleados")
$ws.Range("A1").Value = "CODIGO"
$ws.Range("B1").Value = "APELLIDO"
$ws.Range("C1").Value = "NOMBRE"
$ws.Range("D1").Value = "LEGAJO"
$ws.Range("E1").Value = "SECTOR"
$lo = $ws.ListObjects.Item("empleados")
foreach ($col in $lo.ListColumns) {
    Write-Output $col.Name
}
